# "Preventivo finito + update"
# Fill in the hours-per-phase table (rows 2-8, cols B-G) for each team
# member. Cells that have no hours logged keep the literal dash "-"
# (matching the source workbook's convention) instead of staying blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 - Andrea Favero
$ws.Range("B2").Value = "-"
$ws.Range("C2").Value = "-"
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = 4
$ws.Range("G2").Value = 5

# Row 3 - Eleonora Thiella
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = "-"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = 12
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 8

# Row 4 - Federico Caldart
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = "-"
$ws.Range("D4").Value = 5
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = "-"

# Row 5 - Giovanni Cavallin
$ws.Range("B5").Value = "-"
$ws.Range("C5").Value = "-"
$ws.Range("D5").Value = "-"
$ws.Range("E5").Value = 16
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 7

# Row 6 - Giovanni Dalla Riva
$ws.Range("B6").Value = "-"
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = 6
$ws.Range("F6").Value = 6
$ws.Range("G6").Value = 9

# Row 7 - Lorenzo Menegon
$ws.Range("B7").Value = "-"
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 8
$ws.Range("E7").Value = 8
$ws.Range("F7").Value = 5
$ws.Range("G7").Value = "-"

# Row 8 - Stefano Panozzo
$ws.Range("B8").Value = "-"
$ws.Range("C8").Value = 4
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = 12
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 6

# The staging numbers that used to live in row 11 (5, 10, 20, 80, 20, 35)
# have now been superseded by the real SUM() totals in row 9, so row 11
# is cleared out (its cell styling is left untouched).
$ws.Range("B11:G11").ClearContents() | Out-Null

# Move / resize the stacked bar chart (it was dragged further right and
# enlarged). Left/Top/Width/Height are in points and translate back to the
# underlying column/row anchor offsets when the workbook is saved.
$co = $ws.ChartObjects().Item(1)
$co.Left = 695.5897650098425
$co.Top = 16.87503937007874
$co.Width = 747.5000196850394
$co.Height = 286.35

# Reflect the final selection left by the author (the whole filled table).
$ws.Range("A1:H9").Select() | Out-Null
